# Updated symbol list on Mon Dec 19 05:43:17 UTC 2022 with GitHub Actions
#
# Refreshes the "cryptos" price sheet: most rows just get a refreshed
# Price (column D) reading from the feed, a couple of label cells
# (column E "Worstin24h" suffix) flip, and the BKEXToken / CEJI /
# KickToken trio (rows 41-43) rotates one slot down the ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a bare number (the Price column). These must be
# (re-)entered as Text so they keep behaving like the original inlineStr
# cells instead of turning into numeric cells.
$numericTextCells = @{
    "D2"  = "247.42"
    "D3"  = "21.81"
    "D4"  = "5.492"
    "D5"  = "0.05654"
    "D6"  = "3.381"
    "D7"  = "6.433"
    "D8"  = "0.8014"
    "D9"  = "1.037"
    "D10" = "0.1425"
    "D11" = "0.07241"
    "D12" = "0.03146"
    "D14" = "0.09289"
    "D15" = "0.001660"
    "D16" = "3.224"
    "D17" = "0.04737"
    "D18" = "0.0005857"
    "D19" = "0.006362"
    "D20" = "0.005017"
    "D21" = "0.001049"
    "D23" = "0.0003204"
    "D24" = "4.027"
    "D25" = "2.107"
    "D40" = "0.04083"
    "D41" = "0.006938"
    "D42" = "0.1042"
    "D43" = "0.002973"
    "D44" = "0.009147"
    "D46" = "0.00000000751"
    "D47" = "0.7862"
    "D48" = "0.01288"
    "D49" = "0.00002103"
}

foreach ($addr in $numericTextCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextCells[$addr]
}

# Cells whose new value is already non-numeric text (coin names, links,
# and the composite "rank+name+symbol[+Worstin24h]" labels) - these stay
# text automatically, no special formatting required.
$textCells = @{
    "E18" = "17OneONE"

    "B41" = "KickToken"
    "C41" = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
    "E41" = "40KickTokenKICK"

    "B42" = "BKEXToken"
    "C42" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "E42" = "41BKEXTokenBKK"

    "B43" = "CEJI"
    "C43" = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
    "E43" = "42CEJICEJI"

    "E48" = "47BOLOBOLOWorstin24h"
}

foreach ($addr in $textCells.Keys) {
    $ws.Range($addr).Value = $textCells[$addr]
}
